# Implementation.docx edit ("Add files via upload")
#
# The only semantically meaningful change in the target diff is the
# extension of the "Example: higher layer Resource ..." paragraph with
# two extra sentences:
#
#   "... (interaction context "dialogs")."
#   -> "... (interaction context "dialogs"). Resolution may propagate to
#      other peers (content aware addressing dataflow routes dispatch:
#      P2P resources address encodings, matching forms models requests).
#      Nested interactions."
#
# (The rest of the diff against styles.xml is just several more verbatim
# repeats of the already-repeated built-in Normal/Heading1-6/Title and
# Subtitle style definitions that were sitting in the document before this
# edit. styleId is an identity key in the Word object model - Styles.Add /
# any other automation call against an existing styleId resolves to the
# existing style rather than minting a second one with the same id, so
# that duplication is not something Word's UI/API can (re)produce; it is
# left untouched here.)

$d = $word.ActiveDocument

$old = "peer to “ask” for form elements to be populated (interaction context “dialogs”)."
$new = "peer to “ask” for form elements to be populated (interaction context “dialogs”). Resolution may propagate to other peers (content aware addressing dataflow routes dispatch: P2P resources address encodings, matching forms models requests). Nested interactions."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Find/Replace did not locate the target paragraph text"
}

Write-Output "Paragraph extended: $found"
